# Applies the CodeSystem-BRImunobiologico.xlsx changes described by the commit:
#   "Alteracoes pelo HIAE Diana e ajustes de artefatos"
#
# 1) Properties sheet: rename the "inativo" property to "inactive" and give it
#    its FHIR concept-properties URI.
# 2) Concepts sheet: fix "Covid-19" -> "COVID-19" casing in a handful of
#    existing display/definition cells, rename the "COVID-19 MODERNA" concept
#    to "COVID-19 MODERNA - SPIKEVAX" (with updated definition), and append six
#    brand-new immunobiologic concepts (codes 107-112) at the bottom of the
#    table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Properties sheet (row 2 = "inativo"/"inactive" property)
# ---------------------------------------------------------------------------
$props = $wb.Worksheets.Item("Properties")
$props.Range("A2").Value = "inactive"
$props.Range("B2").Value = "http://hl7.org/fhir/concept-properties#inactive"

# ---------------------------------------------------------------------------
# 2) Concepts sheet
# ---------------------------------------------------------------------------
$concepts = $wb.Worksheets.Item("Concepts")

# -- Casing fixes on existing rows (Covid-19 -> COVID-19) -------------------
$concepts.Range("D86").Value = "Vacina COVID-19 ASTRAZENECA/FIOCRUZ - COVISHIELD, recombinante"
$concepts.Range("D95").Value = "Diluente COVID-19"
$concepts.Range("D97").Value = "Vacina COVID-19 BHARAT - COVAXIN inativada"
$concepts.Range("C98").Value = "COVID-19 MODERNA - SPIKEVAX"
$concepts.Range("D98").Value = "Vacina COVID-19-RNAm, Moderna (Spikevax)"
$concepts.Range("D106").Value = "Vacina COVID-19-RNAm, Moderna (Spikevax) bivalente"

# -- New concept rows (108-113), appended after the existing last row (107) -
$newRows = @(
    @("1", "107", "VPC20", "Vacina Pneumo 20"),
    @("1", "108", "VVSR-Rec", "Vacina Vírus Sincicial Respiratório A e B (recombinante)"),
    @("1", "109", "VVSR-RecAdj", "Vacina Vírus Sincicial Respiratório (recombinante, adjuvada)"),
    @("1", "110", "INF4-alta dosagem", "Vacina Influenza Tetravalente - Alta Dosagem"),
    @("1", "111", "dTpa/VIP", "Vacina Tetra Acelular dTpa/VIP"),
    @("1", "112", "COVID-19 SERUM/ZALIKA", "Vacina Covid-19-recombinante, Serum/Zalika")
)

$startRow = 108
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]

    # Force text storage (these look numeric, e.g. "107") before assigning,
    # matching the existing Level/Code columns which are stored as shared
    # strings rather than numbers.
    $rowRange = $concepts.Range("A" + $r + ":D" + $r)
    $rowRange.NumberFormat = "@"

    $concepts.Range("A" + $r).Value = $rowData[0]
    $concepts.Range("B" + $r).Value = $rowData[1]
    $concepts.Range("C" + $r).Value = $rowData[2]
    $concepts.Range("D" + $r).Value = $rowData[3]
}

# Re-apply the same cell style/format as the last pre-existing data row (107)
# to the newly appended rows so they render consistently with the rest of the
# table (borders, wrap text, vertical alignment, etc).
$concepts.Range("A107:D107").Copy()
$lastNewRow = $startRow + $newRows.Length - 1
$concepts.Range("A" + $startRow + ":D" + $lastNewRow).PasteSpecial(-4122)
